# Continuação do guia de arquitetura
# Update several "D" column cells on the single worksheet to reflect
# progress on the architecture guide (several placeholder "Exemplo:"
# cells are replaced with either concrete statuses or removed to
# "Em desenvolvimento", etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D3: Interfaces de integração -> "Chat " (was "Exemplo: Web Service RESTAPI")
$ws.Range("D3").Value = "Chat "

# D4: Implementação dos recursos de log -> "Em desenvolvimento"
$ws.Range("D4").Value = "Em desenvolvimento"

# D6: Classe de comunicação com os sensores -> "Em desenvolvimento"
$ws.Range("D6").Value = "Em desenvolvimento"

# D10: Configuração da IDE de deploy automatizado -> "Em desenvolvimento"
$ws.Range("D10").Value = "Em desenvolvimento"

# D11: Definição do modelo de autenticação -> "Em desenvolvimento"
$ws.Range("D11").Value = "Em desenvolvimento"

# D15: Servidor de Aplicações, SDKs, Bibliotecas -> plain "Java 11"
# (replaces rich-text "Java" + red ", ? IIS ? Tomcat ?")
$ws.Range("D15").Value = "Java 11"
$ws.Range("D15").Font.Color = 0
$ws.Range("D15").Font.ColorIndex = 1

# D20: Definição de como o software será documentado -> " Word, Excell, Canva"
# (was "Notion, Word, Excell, Canva")
$ws.Range("D20").Value = " Word, Excell, Canva"

# D25: Processo e ferramenta para realização dos Testes -> "Em desenvolvimento"
# (replaces rich-text "? TesteCase " + underlined red "+ Ferramenta ? Jmeter")
$ws.Range("D25").Value = "Em desenvolvimento"
$ws.Range("D25").Font.Color = 0
$ws.Range("D25").Font.ColorIndex = 1
$ws.Range("D25").Font.Underline = $false

# Update the view: scroll so row 7 is the top-most visible row, and
# make D20 the active/selected cell.
$ws.Range("D20").Select()
$excel.ActiveWindow.ScrollRow = 7
